$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new week's worth of data (one "Primera" + one "Segunda" row) is being
# inserted at the top of the Cilantro price table (row 160), pushing the
# existing rows 160:265 down to 162:267.
$ws.Rows("160:161").Insert()

# New row 160 - "Primera"
$ws.Range("A160").Value = 7
$ws.Range("B160").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C160").Value = "Ñuble"
$ws.Range("D160").Value = 45096
$ws.Range("E160").Value = 16
$ws.Range("F160").Value = 100112040
$ws.Range("G160").Value = "Cilantro"
$ws.Range("H160").Value = "Sin especificar"
$ws.Range("I160").Value = "Primera"
$ws.Range("J160").Value = 240
$ws.Range("K160").Value = 1000
$ws.Range("L160").Value = 1200
$ws.Range("M160").Value = 1100
$ws.Range("N160").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O160").Value = "Provincia de Diguillín"
$ws.Range("P160").Value = 1100
$ws.Range("Q160").Value = 1
$ws.Range("R160").Value = "Hortaliza"

# New row 161 - "Segunda"
$ws.Range("A161").Value = 7
$ws.Range("B161").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C161").Value = "Ñuble"
$ws.Range("D161").Value = 45096
$ws.Range("E161").Value = 16
$ws.Range("F161").Value = 100112040
$ws.Range("G161").Value = "Cilantro"
$ws.Range("H161").Value = "Sin especificar"
$ws.Range("I161").Value = "Segunda"
$ws.Range("J161").Value = 100
$ws.Range("K161").Value = 800
$ws.Range("L161").Value = 800
$ws.Range("M161").Value = 800
$ws.Range("N161").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O161").Value = "Provincia de Diguillín"
$ws.Range("P161").Value = 800
$ws.Range("Q161").Value = 1
$ws.Range("R161").Value = "Hortaliza"
